# "altera data por excel"
# - Rename header cell A2 from "usuario" to "nome"
# - Underline the "usuario 1" label in A3
# - Move the active cell/selection to A3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: "usuario" -> "nome"
$ws.Range("A2").Value = "nome"

# A3: underline the font ("usuario 1")
$ws.Range("A3").Font.Underline = $true

# Update the selected cell to A3
$ws.Range("A3").Select()
